$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet 1")

# Update the GCS coverage descriptions: replace embedded newlines with literal "<br>" markers.
# Every cell sharing the string must be updated so the shared string table collapses back
# to a single (updated) entry instead of splitting into a duplicate.
$ws.Range("E30:E43").Value = "Supports the GCS if coverage is **Low**<br>Other members: Global South + EU<br>(25-33% of world emissions)"
$ws.Range("E44:E57").Value = "Supports the GCS if coverage is **Mid**<br>Global South + China<br>(56% of world emissions)"
$ws.Range("E58:E71").Value = "Supports the GCS if coverage is **High**<br>Global South + China + EU + various HICs<br>(UK, Japan, Korea, Canada...; 64-72% of emissions)"
$ws.Range("E72:E85").Value = "Supports the GCS if coverage is **High**, **color** variant<br>Global South + China + EU + various HICs<br>+ Distributive effects shown using colors on world map"

# Update numeric mean/CI_low/CI_high values for the "All" rows that changed
$ws.Range("B2").Value = 67.8255122017956
$ws.Range("C2").Value = 66.435391002076
$ws.Range("D2").Value = 69.2156334015152

$ws.Range("B16").Value = 55.3746767090015
$ws.Range("C16").Value = 54.4859318721968
$ws.Range("D16").Value = 56.2634215458062

$ws.Range("B44").Value = 67.1354010141054
$ws.Range("C44").Value = 65.427782749948
$ws.Range("D44").Value = 68.8430192782629

$ws.Range("B58").Value = 68.4899159160604
$ws.Range("C58").Value = 66.8593150984853
$ws.Range("D58").Value = 70.1205167336354

$ws.Range("B72").Value = 61.8798692282585
$ws.Range("C72").Value = 60.1485598973454
$ws.Range("D72").Value = 63.6111785591716
